$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RW05 scenarios (rows 96-106): SmokeTest should be "No"
$ws.Range("D96:D106").Value = "No"

# RW06 scenarios (rows 107-117): SmokeTest should be "Yes"
$ws.Range("D107:D117").Value = "Yes"
